# Append five new patient rows (9-13) to the Patients sheet, matching the
# structure of the existing "Moshe Davidian" sample entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("A9").Value = "Moshe"
$ws.Range("B9").Value = "Davidian"
$ws.Range("C9").Value = "315783522"
$ws.Range("D9").Value = "22"
$ws.Range("E9").Value = "Yes"
$ws.Range("F9").Value = "Yes"
$ws.Range("G9").Value = "Yes"
$ws.Range("H9").Value = "Yes"
$ws.Range("I9").Value = "Yes"
$ws.Range("J9").Value = "10000"
$ws.Range("K9").Value = "40"
$ws.Range("L9").Value = "40"
$ws.Range("M9").Value = "5"
$ws.Range("N9").Value = "40"
$ws.Range("O9").Value = "20"
$ws.Range("P9").Value = "15"
$ws.Range("Q9").Value = "0.7"
$ws.Range("R9").Value = "80"
$ws.Range("S9").Value = "50"
$ws.Range("T9").Value = "80"
$ws.Range("U9").Value = "The tests are normal and you are a healthy person."

# Row 10
$ws.Range("A10").Value = "Moshe"
$ws.Range("B10").Value = "Davidian"
$ws.Range("C10").Value = "5783522"
$ws.Range("D10").Value = "90"
$ws.Range("E10").Value = "Yes"
$ws.Range("F10").Value = "Yes"
$ws.Range("G10").Value = "Yes"
$ws.Range("H10").Value = "Yes"
$ws.Range("I10").Value = "Yes"
$ws.Range("J10").Value = "10000"
$ws.Range("K10").Value = "40"
$ws.Range("L10").Value = "40"
$ws.Range("M10").Value = "5"
$ws.Range("N10").Value = "40"
$ws.Range("O10").Value = "20"
$ws.Range("P10").Value = "15"
$ws.Range("Q10").Value = "0.7"
$ws.Range("R10").Value = "80"
$ws.Range("S10").Value = "50"
$ws.Range("T10").Value = "80"

# Row 11
$ws.Range("A11").Value = "Moshe"
$ws.Range("B11").Value = "Davidian"
$ws.Range("C11").Value = "315783522"
$ws.Range("D11").Value = "90"
$ws.Range("E11").Value = "Yes"
$ws.Range("F11").Value = "Yes"
$ws.Range("G11").Value = "Yes"
$ws.Range("H11").Value = "Yes"
$ws.Range("I11").Value = "Yes"
$ws.Range("J11").Value = "10000"
$ws.Range("K11").Value = "40"
$ws.Range("L11").Value = "40"
$ws.Range("M11").Value = "5"
$ws.Range("N11").Value = "40"
$ws.Range("O11").Value = "20"
$ws.Range("P11").Value = "15"
$ws.Range("Q11").Value = "0.7"
$ws.Range("R11").Value = "80"
$ws.Range("S11").Value = "50"
$ws.Range("T11").Value = "80"
$ws.Range("U11").Value = "The tests are normal and you are a healthy person."

# Row 12
$ws.Range("A12").Value = "Moshe"
$ws.Range("B12").Value = "Davidian"
$ws.Range("E12").Value = "No"
$ws.Range("F12").Value = "No"
$ws.Range("G12").Value = "No"
$ws.Range("H12").Value = "No"
$ws.Range("I12").Value = "No"
$ws.Range("J12").Value = "10000"
$ws.Range("K12").Value = "40"
$ws.Range("L12").Value = "40"
$ws.Range("M12").Value = "5"
$ws.Range("N12").Value = "40"
$ws.Range("O12").Value = "20"
$ws.Range("P12").Value = "15"
$ws.Range("Q12").Value = "0.7"
$ws.Range("R12").Value = "80"
$ws.Range("S12").Value = "50"
$ws.Range("T12").Value = "80"

# Row 13
$ws.Range("A13").Value = "Moshe"
$ws.Range("B13").Value = "Davidian"
$ws.Range("C13").Value = "315783522"
$ws.Range("D13").Value = "22"
$ws.Range("E13").Value = "Yes"
$ws.Range("F13").Value = "No"
$ws.Range("G13").Value = "No"
$ws.Range("H13").Value = "No"
$ws.Range("I13").Value = "No"
$ws.Range("J13").Value = "10000"
$ws.Range("K13").Value = "40"
$ws.Range("L13").Value = "40"
$ws.Range("M13").Value = "5"
$ws.Range("N13").Value = "40"
$ws.Range("O13").Value = "20"
$ws.Range("P13").Value = "15"
$ws.Range("Q13").Value = "0.7"
$ws.Range("R13").Value = "80"
$ws.Range("S13").Value = "50"
$ws.Range("T13").Value = "80"
$ws.Range("U13").Value = "The tests are normal and you are a healthy person."
